# Consolidate the 2017 ETNP MS/MS worksheets:
#  - remove the "Cyano peps" sheet (its single value now lives in "DB peps")
#  - rename "prok only" to "DB peps"
#  - leave it as the active sheet, selected at A3

$wb = $excel.ActiveWorkbook

$dbPeps = $wb.Worksheets("prok only")
$dbPeps.Name = "DB peps"

$wb.Worksheets("Cyano peps").Delete()

$dbPeps.Activate()
$dbPeps.Range("A3").Select()
